# Updated qc flag format
#
# Normalises the "flag_values"/"flag_meanings" text for several qc flag
# variables: replace the old "0b,1b, 2b[, 3b]" style flag_values strings
# with plain "0, 1, 2[, 3]" text, and replace the newline-separated
# flag_meanings strings with single (space separated) lines.
#
# NOTE: the original spreadsheet was not edited 100% consistently - some
# occurrences ended up with a single space between words/numbers, others
# with a double space. Cells are touched in sheet (row) order so the
# resulting shared-string table is built up the same way it was originally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$longMeaning = "bad_data good_data _good_for_reasearch suspect_data_good_for_general_use suspect_data_requires_further_checking_but_may_be_ok_for_general_use"

$ws.Range("C636").Value2 = "0, 1, 2"
$ws.Range("C637").Value2 = "bad_data good_data suspect_data"

$ws.Range("C644").Value2 = "0, 1, 2"
$ws.Range("C645").Value2 = "bad_data good_data  suspect_data"

$ws.Range("C652").Value2 = "0, 1, 2"
$ws.Range("C653").Value2 = "bad_data good_data  suspect_data"

$ws.Range("C660").Value2 = "0, 1,  2"
$ws.Range("C661").Value2 = "bad_data good_data suspect_data"

$ws.Range("C668").Value2 = "0, 1,  2"
$ws.Range("C669").Value2 = "bad_data good_data suspect_data"

$ws.Range("C676").Value2 = "0, 1,  2"
$ws.Range("C677").Value2 = "bad_data good_data suspect_data"

$ws.Range("C684").Value2 = "0, 1,  2"
$ws.Range("C685").Value2 = "bad_data good_data suspect_data"

$ws.Range("C692").Value2 = "0, 1, 2"
$ws.Range("C693").Value2 = "bad_data good_data suspect_data"

$ws.Range("C700").Value2 = "0, 1,  2, 3"
$ws.Range("C701").Value2 = $longMeaning

$ws.Range("C708").Value2 = "0, 1,  2, 3"
$ws.Range("C709").Value2 = $longMeaning

$ws.Range("C716").Value2 = "0, 1,  2, 3"
$ws.Range("C717").Value2 = $longMeaning

# --- reflect the scroll/selection position recorded in the saved file ------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 685
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C692:C693").Select()
